$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Introduce the new shared strings in the same order they first appear in
# the saved file ("avg" at index 11, then "step/dist" at index 12).
$ws.Range("F6").Value = "avg"
$ws.Range("D2").Value = "step/dist"

# --- FORWARD block (rows 1-7): step/dist column D and avg in F6/F7 ---
# D3:D5 step/dist formulas (shared formula, relative per row)
$ws.Range("D3:D5").Formula = "=A3/B3"

# average formula for FORWARD+BACKWARD combined
$ws.Range("F7").Formula = "=AVERAGE(D3:D5,D9:D11)"

# --- BACKWARD block (rows 8-11): step/dist column D ---
$ws.Range("D9").Formula = "=A9/B9"
$ws.Range("D10:D11").Formula = "=A10/B10"

# --- STRAFE RIGHT block (rows 13-17): add step/dist column J and avg in L18/L19 ---
$ws.Range("J14").Value = "step/dist"
$ws.Range("J16").Formula = "=A16/B16"
$ws.Range("J17").Formula = "=A17/B17"

$ws.Range("L18").Value = "avg"
$ws.Range("L19").Formula = "=AVERAGE(J16:J17,J22:J23)"

# --- STRAFE LEFT block (rows 19-23): step/dist column J ---
$ws.Range("J22").Formula = "=A22/B22"
$ws.Range("J23").Formula = "=A23/B23"

# Update the active selection to I7, matching the saved workbook view
$ws.Range("I7").Select()
